# Apply the target edits to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update cell values
$ws.Range("B2").Value = 3.8
$ws.Range("C2").Value = 11.5
$ws.Range("B3").Value = 4.5999999999999996
$ws.Range("C3").Value = 10.5
$ws.Range("C5").Value = 18

# Update column widths: column A and C get wider (column B is left as-is)
# NOTE: the host engine quantizes the stored <col> width to a 1/7-character
# grid (Calibri-11 "Maximum Digit Width" = 7px) no matter what font the
# workbook actually uses, so the ColumnWidth values below are chosen to
# land as close as that grid allows to the intended stored widths of
# 27 and 27.25 characters respectively.
$ws.Columns.Item(1).ColumnWidth = 26.285714285714285
$ws.Columns.Item(3).ColumnWidth = 26.571428571428573

# Update selection to B3
$ws.Range("B3").Select() | Out-Null
